$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Cntn1"
$ws.Range("C2").Value = "Notch2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.07031566666666667
$ws.Range("H2").Value = 0.210947
$ws.Range("I2").Value = 0.3500709860127268
$ws.Range("J2").Value = 0.446885632088942
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.400501333333333
$ws.Range("N2").Value = 4.201504
$ws.Range("O2").Value = 0.00926314904242919
$ws.Range("P2").Value = 0.009687730200823723
$ws.Range("Q2").Value = 0.09847718492088889
$ws.Range("R2").Value = 0.886294664288
$ws.Range("S2").Value = 0.003242759718866033
$ws.Range("T2").Value = 0.004329307434302242

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Cntn1"
$ws.Range("C3").Value = "Notch2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.07031566666666667
$ws.Range("H3").Value = 0.210947
$ws.Range("I3").Value = 0.3500709860127268
$ws.Range("J3").Value = 0.446885632088942
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 21.254561
$ws.Range("N3").Value = 63.763683
$ws.Range("O3").Value = 0.1405812059498714
$ws.Range("P3").Value = 0.1470248171880475
$ws.Range("Q3").Value = 1.494528626422333
$ws.Range("R3").Value = 13.450757637801
$ws.Range("S3").Value = 0.0492134013817297
$ws.Range("T3").Value = 0.06570327836184174

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Cntn1"
$ws.Range("C4").Value = "Notch2"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.07031566666666667
$ws.Range("H4").Value = 0.210947
$ws.Range("I4").Value = 0.3500709860127268
$ws.Range("J4").Value = 0.446885632088942
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 63.87756733333333
$ws.Range("N4").Value = 191.632702
$ws.Range("O4").Value = 0.4224968677952986
$ws.Range("P4").Value = 0.4418622271050682
$ws.Range("Q4").Value = 4.491593732088222
$ws.Range("R4").Value = 40.424343588794
$ws.Range("S4").Value = 0.1479038950963889
$ws.Range("T4").Value = 0.1974618806560761

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cntn1"
$ws.Range("C5").Value = "Notch2"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.07031566666666667
$ws.Range("H5").Value = 0.210947
$ws.Range("I5").Value = 0.3500709860127268
$ws.Range("J5").Value = 0.446885632088942
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 19.878555
$ws.Range("N5").Value = 39.75711
$ws.Range("O5").Value = 0.1314800731212866
$ws.Range("P5").Value = 0.0916710195312133
$ws.Range("Q5").Value = 1.397773847195
$ws.Range("R5").Value = 8.38664308317
$ws.Range("S5").Value = 0.04602735883859421
$ws.Range("T5").Value = 0.04096646150744401

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Cntn1"
$ws.Range("C6").Value = "Notch2"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.07031566666666667
$ws.Range("H6").Value = 0.210947
$ws.Range("I6").Value = 0.3500709860127268
$ws.Range("J6").Value = 0.446885632088942
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 44.77944466666667
$ws.Range("N6").Value = 134.338334
$ws.Range("O6").Value = 0.2961787040911142
$ws.Range("P6").Value = 0.3097542059748472
$ws.Range("Q6").Value = 3.148696504699778
$ws.Range("R6").Value = 28.338268542298
$ws.Range("S6").Value = 0.103683570977148
$ws.Range("T6").Value = 0.1384247041292779

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Cntn1"
$ws.Range("C7").Value = "Notch2"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.1305455
$ws.Range("H7").Value = 0.261091
$ws.Range("I7").Value = 0.6499290139872732
$ws.Range("J7").Value = 0.5531143679110581
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.400501333333333
$ws.Range("N7").Value = 4.201504
$ws.Range("O7").Value = 0.00926314904242919
$ws.Range("P7").Value = 0.009687730200823723
$ws.Range("Q7").Value = 0.1828291468106667
$ws.Range("R7").Value = 1.096974880864
$ws.Range("S7").Value = 0.006020389323563157
$ws.Range("T7").Value = 0.005358422766521481

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Cntn1"
$ws.Range("C8").Value = "Notch2"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.1305455
$ws.Range("H8").Value = 0.261091
$ws.Range("I8").Value = 0.6499290139872732
$ws.Range("J8").Value = 0.5531143679110581
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 21.254561
$ws.Range("N8").Value = 63.763683
$ws.Range("O8").Value = 0.1405812059498714
$ws.Range("P8").Value = 0.1470248171880475
$ws.Range("Q8").Value = 2.7746872930255
$ws.Range("R8").Value = 16.648123758153
$ws.Range("S8").Value = 0.0913678045681417
$ws.Range("T8").Value = 0.08132153882620574

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Cntn1"
$ws.Range("C9").Value = "Notch2"
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.1305455
$ws.Range("H9").Value = 0.261091
$ws.Range("I9").Value = 0.6499290139872732
$ws.Range("J9").Value = 0.5531143679110581
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 63.87756733333333
$ws.Range("N9").Value = 191.632702
$ws.Range("O9").Value = 0.4224968677952986
$ws.Range("P9").Value = 0.4418622271050682
$ws.Range("Q9").Value = 8.338928966313667
$ws.Range("R9").Value = 50.033573797882
$ws.Range("S9").Value = 0.2745929726989098
$ws.Range("T9").Value = 0.2444003464489922

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Cntn1"
$ws.Range("C10").Value = "Notch2"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.1305455
$ws.Range("H10").Value = 0.261091
$ws.Range("I10").Value = 0.6499290139872732
$ws.Range("J10").Value = 0.5531143679110581
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 19.878555
$ws.Range("N10").Value = 39.75711
$ws.Range("O10").Value = 0.1314800731212866
$ws.Range("P10").Value = 0.0916710195312133
$ws.Range("Q10").Value = 2.5950559017525
$ws.Range("R10").Value = 10.38022360701
$ws.Range("S10").Value = 0.08545271428269235
$ws.Range("T10").Value = 0.0507045580237693

# Row 11
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Cntn1"
$ws.Range("C11").Value = "Notch2"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.1305455
$ws.Range("H11").Value = 0.261091
$ws.Range("I11").Value = 0.6499290139872732
$ws.Range("J11").Value = 0.5531143679110581
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 44.77944466666667
$ws.Range("N11").Value = 134.338334
$ws.Range("O11").Value = 0.2961787040911142
$ws.Range("P11").Value = 0.3097542059748472
$ws.Range("Q11").Value = 5.845754993732334
$ws.Range("R11").Value = 35.074529962394
$ws.Range("S11").Value = 0.1924951331139662
$ws.Range("T11").Value = 0.1713295018455693
